# Update countries & provincias Spain
# Refresh the COVID-19 country statistics table (sheet "Pais") with the
# latest snapshot, including the re-ranking of a few countries whose
# case counts crossed each other, and bump the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: last-updated timestamp moves from 16:45 to 18:02
$ws.Range("A1").Value = 'Datos actualizados a 25 de Septiembre de 2020 a las 18:02'

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7195797
$ws.Range("C4").Value = 10326
$ws.Range("D4").Value = 4440476
$ws.Range("E4").Value = 2547546
$ws.Range("G4").Value = 237
$ws.Range("H4").Value = 207775

# Row 23 - Italia
$ws.Range("B23").Value = 306235
$ws.Range("C23").Value = 1912
$ws.Range("D23").Value = 222716
$ws.Range("E23").Value = 47718
$ws.Range("G23").Value = 20
$ws.Range("H23").Value = 35801

# Row 29 - Canada
$ws.Range("B29").Value = 149503
$ws.Range("C29").Value = 409
$ws.Range("D29").Value = 128990
$ws.Range("E29").Value = 11263
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 9250

# Row 34 - Republica Dominicana
$ws.Range("B34").Value = 110597
$ws.Range("C34").Value = 475
$ws.Range("D34").Value = 84610
$ws.Range("E34").Value = 23900
$ws.Range("G34").Value = 11
$ws.Range("H34").Value = 2087

# Row 45 - Guatemala
$ws.Range("B45").Value = 88878
$ws.Range("C45").Value = 945
$ws.Range("D45").Value = 77750
$ws.Range("E45").Value = 7942
$ws.Range("G45").Value = 16
$ws.Range("H45").Value = 3186

# Row 59 - Singapur
$ws.Range("D59").Value = 57341
$ws.Range("E59").Value = 297

# Row 87 - Grecia
$ws.Range("B87").Value = 16913
$ws.Range("C87").Value = 286
$ws.Range("E87").Value = 6555
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = 369

# Row 95 - Albania
$ws.Range("B95").Value = 13045
$ws.Range("C95").Value = 124
$ws.Range("D95").Value = 7309
$ws.Range("E95").Value = 5363
$ws.Range("G95").Value = 3
$ws.Range("H95").Value = 373

# Rows 100-102 - Montenegro overtakes Maldivas and Guayana Francesa
$ws.Range("A100").Value = 'Montenegro'
$ws.Range("B100").Value = 9962
$ws.Range("C100").Value = 245
$ws.Range("D100").Value = 6177
$ws.Range("E100").Value = 3630
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 155

$ws.Range("A101").Value = 'Maldivas'
$ws.Range("B101").Value = 9939
$ws.Range("D101").Value = 8597
$ws.Range("E101").Value = 1308
$ws.Range("H101").Value = 34

$ws.Range("A102").Value = 'Guayana Francesa'
$ws.Range("B102").Value = 9790
$ws.Range("D102").Value = 9456
$ws.Range("E102").Value = 269
$ws.Range("H102").Value = 65

# Row 108 - Luxemburgo
$ws.Range("B108").Value = 8233
$ws.Range("C108").Value = 75
$ws.Range("D108").Value = 6976
$ws.Range("E108").Value = 1133

# Row 114 - Jordania
$ws.Range("B114").Value = 7211
$ws.Range("C114").Value = 620
$ws.Range("D114").Value = 4035
$ws.Range("E114").Value = 3137
$ws.Range("G114").Value = 3
$ws.Range("H114").Value = 39

# Row 153 - Yemen
$ws.Range("D153").Value = 1255
$ws.Range("E153").Value = 187
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 587

# Row 184 - Isla de Man
$ws.Range("D184").Value = 314
$ws.Range("E184").Value = 2

# Rows 215-216 - Montserrat overtakes Islas Malvinas
$ws.Range("A215").Value = 'Montserrat'
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

$ws.Range("A216").Value = 'Islas Malvinas'
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
